# SECTOR_holdings.xlsx update:
#   - bump the "Model holdings provided as of ..." date in the confidential
#     disclaimer from 2021-03-25 to 2021-03-26
#   - refresh the Weight / Percent Change figures in D2:E6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected, so unlock it before writing, then restore
# protection afterwards.
$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect()
}

$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."
# Writing a multi-line value auto-expands the row height; put it back to
# the sheet's default (the source file has no explicit row height here).
$ws.Rows(9).AutoFit()

$ws.Range("D2").Value = 0.2522921400305184
$ws.Range("E2").Value = 0.02481995884773647

$ws.Range("D3").Value = 0.25033752338478
$ws.Range("E3").Value = 0.01566193853427866

$ws.Range("D4").Value = 0.2488953892422488
$ws.Range("E4").Value = 0.02542896892873703

$ws.Range("D5").Value = 0.2484749473424527
$ws.Range("E5").Value = -0.01094540976877811

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.01379214444084909

if ($wasProtected) {
    $ws.Protect()
}
